# Update "cryptos" list values (Price and Volume(1h) columns) with latest
# scraped data from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '30.348.19'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.96%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.859.19'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -1.98%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '234.59'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.65%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4737'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.85%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.2743'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -3.11%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06433'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.55%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.832.65'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.75%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07466'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.19%  '
$ws.Cells.Item(12, 5).Value = '  -1.91%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.993'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.97%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '85.55'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -2.93%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6330'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -4.79%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '30.309.19'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.01%  '
$ws.Cells.Item(17, 5).Value = '  +0.06%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '230.49'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.62%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.75'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000007412'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.25%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.099.23'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -4.35%  '
$ws.Cells.Item(22, 5).Value = '  +0.20%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.994'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -5.75%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.999'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -3.54%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.273'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.19%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '166.39'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.16%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.98'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.892'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -3.09%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1048'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +7.76%  '
$ws.Cells.Item(30, 5).Value = '  -0.11%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.151'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -4.66%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.929'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -2.19%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.04944'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -2.14%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.165'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -5.11%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7249'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -3.94%  '
$ws.Cells.Item(36, 5).Value = '  +0.31%  '
$ws.Cells.Item(37, 5).Value = '  -0.37%  '
$ws.Cells.Item(38, 5).Value = '  -0.34%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.649'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.30%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9157'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.13%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.975'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -5.16%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '106.23'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.08%  '
$ws.Cells.Item(43, 5).Value = '  -0.21%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4113'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -3.92%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.580'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -3.54%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.125'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.70%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '60.93'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -5.37%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1199'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -5.66%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.653'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -3.60%  '
$ws.Cells.Item(50, 5).Value = '  -0.73%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.410'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -4.51%  '
